$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 280637.47
$ws.Range("I17").Value = 1663
$ws.Range("J17").Value = 288608.16
$ws.Range("K17").Value = 4989
$ws.Range("L17").Value = 865824.48
$ws.Range("M17").Value = -4821
$ws.Range("N17").Value = -866160.48
$ws.Range("H51").Value = 64065148
$ws.Range("J51").Value = 3080
$ws.Range("L51").Value = 3080
$ws.Range("N51").Value = -4048
$ws.Range("H53").Value = 627.4545000000001
$ws.Range("I53").Value = 100.5
$ws.Range("J53").Value = 928.5714
$ws.Range("K53").Value = 100.5
$ws.Range("L53").Value = 928.5714
$ws.Range("M53").Value = 536.5
$ws.Range("N53").Value = -2202.5714
$ws.Range("H132").Value = 28580654
$ws.Range("I132").Value = 35719190
$ws.Range("J132").Value = 26508.428
$ws.Range("K132").Value = 107157570
$ws.Range("L132").Value = 79525.284
$ws.Range("M132").Value = -107155040
$ws.Range("N132").Value = -84585.284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7408.702
$ws.Range("I32").Value = 7371.049
$ws.Range("K32").Value = 7371.049
$ws.Range("M32").Value = -7084.049
$ws.Range("H61").Value = 1856.3429
$ws.Range("I61").Value = 1866.7646
$ws.Range("J61").Value = 1502
$ws.Range("K61").Value = 1866.7646
$ws.Range("L61").Value = 1502
$ws.Range("M61").Value = -1654.7646
$ws.Range("N61").Value = -1926
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 1922.2142
$ws.Range("I122").Value = 1839.3077
$ws.Range("K122").Value = 5517.9231
$ws.Range("M122").Value = -3067.9231
$ws.Range("H136").Value = 1856.3429
$ws.Range("I136").Value = 1866.7646
$ws.Range("J136").Value = 1502
$ws.Range("K136").Value = 5600.293799999999
$ws.Range("L136").Value = 4506
$ws.Range("M136").Value = -3050.293799999999
$ws.Range("N136").Value = -9606
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2668.077
$ws.Range("I20").Value = 2784.5386
$ws.Range("K20").Value = 2784.5386
$ws.Range("M20").Value = -2537.5386
$ws.Range("H124").Value = 139999
$ws.Range("J124").Value = 139999
$ws.Range("L124").Value = 139999
$ws.Range("N124").Value = -149819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 35000
$ws.Range("J55").Value = 35000
$ws.Range("L55").Value = 35000
$ws.Range("N55").Value = -35630
$ws.Range("H58").Value = 2468.762
$ws.Range("I58").Value = 2441.1667
$ws.Range("K58").Value = 2441.1667
$ws.Range("M58").Value = -2238.1667
$ws.Range("H107").Value = 3854.3333
$ws.Range("I107").Value = 3061.75
$ws.Range("J107").Value = 4488.4
$ws.Range("K107").Value = 3061.75
$ws.Range("L107").Value = 4488.4
$ws.Range("M107").Value = -1141.75
$ws.Range("N107").Value = -8328.4
$ws.Range("H135").Value = 113499.75
$ws.Range("I135").Value = 104000
$ws.Range("K135").Value = 104000
$ws.Range("M135").Value = -98930
$ws.Range("H136").Value = 2468.762
$ws.Range("I136").Value = 2441.1667
$ws.Range("K136").Value = 7323.500100000001
$ws.Range("M136").Value = -4773.500100000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5952390
$ws.Range("I4").Value = 7405182.5
$ws.Range("K4").Value = 22215547.5
$ws.Range("M4").Value = -22215435.5
$ws.Range("H130").Value = 8999.5
$ws.Range("I130").Value = 8999.5
$ws.Range("K130").Value = 26998.5
$ws.Range("M130").Value = -21978.5
$ws.Range("H131").Value = 1705.4
$ws.Range("I131").Value = 1230
$ws.Range("J131").Value = 1739.3572
$ws.Range("K131").Value = 3690
$ws.Range("L131").Value = 5218.071599999999
$ws.Range("M131").Value = 1350
$ws.Range("N131").Value = -15298.0716
$ws.Range("H140").Value = 3656.7
$ws.Range("I140").Value = 2292.182
$ws.Range("K140").Value = 6876.545999999999
$ws.Range("M140").Value = -1696.545999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3340.6316
$ws.Range("I80").Value = 3189.3333
$ws.Range("K80").Value = 3189.3333
$ws.Range("M80").Value = -2191.3333
$ws.Range("H83").Value = 3340.6316
$ws.Range("I83").Value = 3189.3333
$ws.Range("K83").Value = 15946.6665
$ws.Range("M83").Value = -10954.6665
$ws.Range("H122").Value = 1666.2059
$ws.Range("I122").Value = 1559.1724
$ws.Range("K122").Value = 4677.5172
$ws.Range("M122").Value = -2227.5172
$ws.Range("H132").Value = 2672.2144
$ws.Range("I132").Value = 2672.2144
$ws.Range("K132").Value = 8016.6432
$ws.Range("M132").Value = -5486.6432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1400.9375
$ws.Range("I16").Value = 1222.6428
$ws.Range("K16").Value = 1222.6428
$ws.Range("M16").Value = -1052.6428
$ws.Range("H25").Value = 34900
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H93").Value = 5120.5
$ws.Range("I93").Value = 5560.6665
$ws.Range("J93").Value = 3800
$ws.Range("K93").Value = 5560.6665
$ws.Range("L93").Value = 3800
$ws.Range("M93").Value = -4312.6665
$ws.Range("N93").Value = -6296
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H125").Value = 71877.57000000001
$ws.Range("J125").Value = 71877.57000000001
$ws.Range("L125").Value = 71877.57000000001
$ws.Range("N125").Value = -81717.57000000001
$ws.Range("H136").Value = 3294.7805
$ws.Range("I136").Value = 2951.2
$ws.Range("K136").Value = 8853.599999999999
$ws.Range("M136").Value = -6303.599999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 26277.5
$ws.Range("J76").Value = 26277.5
$ws.Range("L76").Value = 26277.5
$ws.Range("N76").Value = -26907.5
$ws.Range("H79").Value = 26277.5
$ws.Range("J79").Value = 26277.5
$ws.Range("L79").Value = 26277.5
$ws.Range("N79").Value = -28461.5
$ws.Range("H97").Value = 79999.336
$ws.Range("J97").Value = 79999.336
$ws.Range("L97").Value = 79999.336
$ws.Range("N97").Value = -81981.336
$ws.Range("H122").Value = 3416.973
$ws.Range("I122").Value = 3337.76
$ws.Range("K122").Value = 10013.28
$ws.Range("M122").Value = -7563.280000000001
$ws.Range("H125").Value = 500357.5
$ws.Range("J125").Value = 500357.5
$ws.Range("L125").Value = 500357.5
